# AHLE scenario parameters POULTRY.xlsx
# "new small ruminant results" - add AHLE Vil ind summary block to the
# About sheet, add a row-6 total (SUM) on Sheet1, and update the saved
# view/selection state to match the author's final screen position.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("About")

# ---------------------------------------------------------------------
# Sheet1: row 6 gains a running total in column U (=SUM(U7:U12)), plus
# two formatted-but-empty cells (C6, L6) picking up the same "integer on
# light-blue fill" style already used for A6 in that row.
# ---------------------------------------------------------------------
$ws1.Range("A6").Copy()
$ws1.Range("C6").PasteSpecial(-4122) # xlPasteFormats
$ws1.Range("C6").NumberFormat = "0"

$ws1.Range("A6").Copy()
$ws1.Range("L6").PasteSpecial(-4122) # xlPasteFormats
$ws1.Range("L6").NumberFormat = "0"

# New column L needs an explicit width so it shows up in <cols>.
$ws1.Columns.Item(12).ColumnWidth = 8.3

$ws1.Range("U6").Formula = "=SUM(U7:U12)"

# ---------------------------------------------------------------------
# About sheet: new "AHLE Vil ind" small-ruminant results block (rows
# 30-35), plus a couple of stray formatted-but-empty cells (rows 38-39)
# carrying forward the same style used above them in column F.
# ---------------------------------------------------------------------
$ws2.Range("F17").Copy()
$ws2.Range("F30:F32").PasteSpecial(-4122) # xlPasteFormats
$ws2.Range("F38:F39").PasteSpecial(-4122) # xlPasteFormats

$ws2.Range("E30").Value = "vill ind current"
$ws2.Range("F30").Value = 9154907178.8783302

$ws2.Range("E31").Value = "vil ind juv mort zero"
$ws2.Range("F31").Value = 10954558648.7402

$ws2.Range("E32").Value = "vil ind juv ideal"
$ws2.Range("F32").Value = 12965788832.187099

# Shared strings must come out in the same order as the target file, so
# set G33 ("AHLE") before F33 ("Juv").
$ws2.Range("G33").Value = "AHLE"
$ws2.Range("F33").Value = "Juv"
$ws2.Range("H33").Formula = "=F32-F30"

$ws2.Range("F34").Value = "Juv"
$ws2.Range("G34").Value = "AHLE mort"
$ws2.Range("H34").Formula = "=F31-F30"

$ws2.Range("G35").Value = "prod AHLE"
$ws2.Range("H35").Formula = "=H33-H34"

# New column H needs an explicit width so it shows up in <cols>.
$ws2.Columns.Item(8).ColumnWidth = 10.3

# ---------------------------------------------------------------------
# View / selection state: the workbook was last saved with "About" as
# the active sheet/tab, scrolled so row 23 is at the top, with
# E38:F40 selected; Sheet1 keeps its frozen panes but the active cell
# there moves to U6 (the new total).
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("U6").Select()

$ws2.Activate()
$ws2.Range("B23").Select()
$ws2.Range("E38:F40").Select()
